# Column inversion: swap the Longitude (I) / Latitude (J) columns for the
# coordinate data rows (2-4). The cell that carries the thousands-separator
# number format ("#,##0", style index 3) previously lived on J3 and must
# travel with its value to I3, leaving J3 with the default (General) format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 45.467039999999997
$ws.Range("J2").Value = 8.8957999999999995

# Row 3 (carries the special number format that must move from J3 to I3)
$ws.Range("I3").Value = 45.468019499999997
$ws.Range("J3").Value = 8.8923500999999998
$ws.Range("J3").ClearFormats()
$ws.Range("I3").NumberFormat = "#,##0"

# Row 4
$ws.Range("I4").Value = 45.4653931
$ws.Range("J4").Value = 8.8838428999999994

# Update the active selection to I4 (was D8)
$null = $ws.Range("I4").Select()
